# "fixed error in BOM"
# - Corrected unit prices for C1/C2 series caps (R13/R14: 0.8 -> 0.08)
# - Added missing hyperlinks for the GRM155Z71A105KE01J (H15) and LMK105BJ225KV-F (H16) links
# - Cleared a stray fill formatting override on A19 so it matches the rest of the row
# - Updated the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the unit price typo for rows 13 and 14 (0.8 was a 10x error, should be 0.08) ---
$ws.Range("E13").Value = 0.08
$ws.Range("E14").Value = 0.08

# --- Row 19 (J2 / Front Panel USB 2.0 Header): clear the stray fill so the style matches the rest of the row ---
$ws.Range("A19").Interior.Pattern = -4142

# --- Add the missing hyperlink for H15 (GRM155Z71A105KE01J crystal load cap), text unchanged ---
$ws.Hyperlinks.Add($ws.Range("H15"), "https://www.digikey.com/en/products/detail/murata-electronics/GRM155Z71A105KE01J/16033607") | Out-Null

# --- Add the missing hyperlink for H16 (LMK105BJ225KV-F crystal), with trailing-space text to match the other link cells ---
$ws.Range("H16").Value = "https://www.digikey.com/en/products/detail/taiyo-yuden/LMK105BJ225KV-F/7403747 "
$ws.Hyperlinks.Add($ws.Range("H16"), "https://www.digikey.com/en/products/detail/taiyo-yuden/LMK105BJ225KV-F/7403747") | Out-Null

# --- Move the active selection to C1 ---
$ws.Range("C1").Select() | Out-Null
